$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Treatment Agent" column in the TreatmentTab query (B5) ---
# Replace CONCAT(REPLACE(...)) with a plain REPLACE(...) call.
$oldAgentExpr = "CONCAT(REPLACE(trt.treatment_agent, ';', ', '))"
$newAgentExpr = "REPLACE(trt.treatment_agent, ';', ', ')"

$treatmentQuery = $ws.Range("B5").Value2
$treatmentQuery = $treatmentQuery.Replace($oldAgentExpr, $newAgentExpr)
$ws.Range("B5").Value = $treatmentQuery

# Re-apply the cell's font explicitly (matches the formatting touch-up that
# was present in the authored workbook for this cell).
$ws.Range("B5").Font.Name = "Calibri"
$ws.Range("B5").Font.Size = 12
$ws.Range("B5").WrapText = $true

# --- Update the saved view/selection state ---
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("C5").Select() | Out-Null
